# Agent Commission Setup completed
#
# The "ProducerCode" sheet held two separate producer codes (A1: AG1730A1,
# A2: AG8160A1). They get consolidated into a single updated code in A1,
# the now-empty second row is removed, and A1 is given a "form field" look
# (Text number format, thin box border, wrapped text, wider column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is no longer needed now that the codes are consolidated into A1.
$ws.Rows(2).Delete()

# New consolidated producer code.
$ws.Range("A1").Value2 = "AG6304A44"

# Give A1 a bordered, text-formatted, word-wrapped look.
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").WrapText = $true

# Widen column A to fit the new look.
$ws.Columns("A").ColumnWidth = 26.5

# Leave the selection where the user ended up after the edit.
[void]$ws.Range("C3").Select()
